$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ELBASVIR")
$ws.Columns("B:B").Cut() | Out-Null
$ws.Columns("A:A").Insert() | Out-Null
